$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3: now holds the Archer data (previously in row 4) ---
$ws.Range("A3").Value = "Archer"
$ws.Range("B3").Value = 0.6
$ws.Range("C3").Value = 5
$ws.Range("E3").Value = "Ranged"
$ws.Range("F3").Value = 1.5
$ws.Range("G3").Value = 75
$ws.Range("O3").Value = "Warrior vs Skeleton"
$ws.Range("P3").Value = "Warrior beats skeleton in ~5 seconds w/ half health remaining"

# --- Row 4: now holds the Skeleton data (previously in row 3) ---
$ws.Range("A4").Value = "Skeleton"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 10
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 50

# Clear the leftover cells that no longer apply to row 4
$ws.Range("E4").Clear() | Out-Null
$ws.Range("O4").Clear() | Out-Null
$ws.Range("P4").Clear() | Out-Null

# --- Row 5: Ghoul damage-per-hit tweak ---
$ws.Range("B5").Value = 0.75

# --- Row 6: new Vampire entry ---
$ws.Range("A6").Value = "Vampire"
$ws.Range("B6").Value = 0.5
$ws.Range("C6").Value = 10
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 100

# --- Row 7: new Zombie entry ---
$ws.Range("A7").Value = "Zombie"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 20
$ws.Range("F7").Value = 1.5
$ws.Range("G7").Value = 150

# Second pass for the new ability-name strings, so the shared-string table
# gets Vampire(18), Zombie(19), Life Steal(20), Stun(21), Mummy(22), Vulnerable(23)
$ws.Range("H6").Value = "Life Steal"
$ws.Range("I6").Formula = "=(G6 + D6*0.25*5)/G6"

$ws.Range("E7").Value = "Stun"

# --- Row 8: new Mummy entry ---
$ws.Range("A8").Value = "Mummy"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 15
$ws.Range("E8").Value = "Vulnerable"
$ws.Range("F8").Value = 1.25
$ws.Range("G8").Value = 200

# Re-apply the shared DPS formula on the rows whose inputs changed after the
# formula was first evaluated, so the cached results reflect the final inputs.
$ws.Range("J6").Formula = "=D6*F6/10 +G6*I6/100"
$ws.Range("J7").Formula = "=D7*F7/10 +G7*I7/100"
$ws.Range("J8").Formula = "=D8*F8/10 +G8*I8/100"

$ws.Range("A13").Select() | Out-Null
$wb.Application.CalculateFull() | Out-Null
